$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5987
$ws1.Range("F9").Value = 49
$ws1.Range("F10").Value = 781
$ws1.Range("F11").Value = 358
$ws1.Range("F12").Value = 4543
$ws1.Range("F13").Value = 4543
$ws1.Range("F15").Value = 105
$ws1.Range("F16").Value = 141
$ws1.Range("F20").Value = 7113
$ws1.Range("F21").Value = 242
$ws1.Range("F22").Value = 118
$ws1.Range("F24").Value = 492
$ws1.Range("F26").Value = 6272
$ws1.Range("F30").Value = 6077
$ws1.Range("F33").Value = 107
$ws1.Range("F35").Value = 442
$ws1.Range("F36").Value = 6202
$ws1.Range("F47").Value = 383
$ws1.Range("F48").Value = 2097

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 219
$ws2.Range("F5").Value = 38
$ws2.Range("F6").Value = 109

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5988
$ws4.Range("F4").Value = 5988
$ws4.Range("F8").Value = 219
$ws4.Range("F9").Value = 49
$ws4.Range("F10").Value = 38
$ws4.Range("F11").Value = 358
$ws4.Range("F12").Value = 4543
$ws4.Range("F13").Value = 4543
$ws4.Range("F15").Value = 105
$ws4.Range("F16").Value = 141
$ws4.Range("F20").Value = 7113
$ws4.Range("F21").Value = 242
$ws4.Range("F22").Value = 118
$ws4.Range("F23").Value = 492
$ws4.Range("F25").Value = 109
$ws4.Range("F26").Value = 6272
$ws4.Range("F31").Value = 6077
$ws4.Range("F35").Value = 107
$ws4.Range("F37").Value = 442
$ws4.Range("F38").Value = 6202
$ws4.Range("F48").Value = 383
$ws4.Range("F49").Value = 2097

